$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$counts = @(2,0,75,52,3,116,123,42,124,61,18,69,87,83,100,77,44,5,47,35,28,22,30,27,86,39,29,15,58,24,93,115)
$images = @("flower/flower003.jpg","face/face030.jpg","face/face026.jpg","flower/flower025.jpg","face/face016.jpg","face/face006.jpg","flower/flower024.jpg","face/face022.jpg","flower/flower030.jpg","flower/flower001.jpg","face/face002.jpg","flower/flower031.jpg","flower/flower006.jpg","flower/flower029.jpg","face/face025.jpg","flower/flower028.jpg","face/face008.jpg","face/face019.jpg","face/face028.jpg","flower/flower017.jpg","face/face004.jpg","flower/flower008.jpg","face/face015.jpg","flower/flower005.jpg","face/face018.jpg","flower/flower015.jpg","flower/flower011.jpg","face/face003.jpg","face/face001.jpg","flower/flower010.jpg","face/face029.jpg","flower/flower023.jpg")
$words = @("lassen","wehen","hören","duschen","nullen","fließen","spenden","kranken","dienen","münzen","tollen","heißen","wählen","binden","proben","trotzen","achten","lächeln","frischen","zögern","passen","ändern","wecken","quellen","parken","holen","spüren","mögen","herrschen","kriegen","prüfen","planen")
$cats = @("flower","face","face","flower","face","face","flower","face","flower","flower","face","flower","flower","flower","face","flower","face","face","face","flower","face","flower","face","flower","face","flower","flower","face","face","flower","face","flower")

for ($i = 0; $i -lt $images.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $images[$i]
}
for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $words[$i]
}
for ($i = 0; $i -lt $cats.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $cats[$i]
}
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}
